$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88-59=29"
$t.Cell(1, 2).Range.Text = "41-2=39"
$t.Cell(1, 3).Range.Text = "50-17=33"
$t.Cell(1, 4).Range.Text = "94-76=18"
$t.Cell(1, 5).Range.Text = "68+16=84"
$t.Cell(2, 1).Range.Text = "96-77=19"
$t.Cell(2, 2).Range.Text = "56-29=27"
$t.Cell(2, 3).Range.Text = "83-18=65"
$t.Cell(2, 4).Range.Text = "26+6=32"
$t.Cell(2, 5).Range.Text = "4+57=61"
$t.Cell(3, 1).Range.Text = "77-29=48"
$t.Cell(3, 2).Range.Text = "5+67=72"
$t.Cell(3, 3).Range.Text = "8+56=64"
$t.Cell(3, 4).Range.Text = "39+9=48"
$t.Cell(3, 5).Range.Text = "61-6=55"
$t.Cell(4, 1).Range.Text = "90-41=49"
$t.Cell(4, 2).Range.Text = "7+76=83"
$t.Cell(4, 3).Range.Text = "84-7=77"
$t.Cell(4, 4).Range.Text = "60-55=5"
$t.Cell(4, 5).Range.Text = "58+27=85"
$t.Cell(5, 1).Range.Text = "54+27=81"
$t.Cell(5, 2).Range.Text = "35+19=54"
$t.Cell(5, 3).Range.Text = "9+37=46"
$t.Cell(5, 4).Range.Text = "59+22=81"
$t.Cell(5, 5).Range.Text = "94-16=78"
$t.Cell(6, 1).Range.Text = "26+46=72"
$t.Cell(6, 2).Range.Text = "78+19=97"
$t.Cell(6, 3).Range.Text = "19+64=83"
$t.Cell(6, 4).Range.Text = "55+36=91"
$t.Cell(6, 5).Range.Text = "70-45=25"
$t.Cell(7, 1).Range.Text = "92-17=75"
$t.Cell(7, 2).Range.Text = "52+9=61"
$t.Cell(7, 3).Range.Text = "77+6=83"
$t.Cell(7, 4).Range.Text = "93-17=76"
$t.Cell(7, 5).Range.Text = "52-26=26"
$t.Cell(8, 1).Range.Text = "75+9=84"
$t.Cell(8, 2).Range.Text = "68+19=87"
$t.Cell(8, 3).Range.Text = "40-27=13"
$t.Cell(8, 4).Range.Text = "55-19=36"
$t.Cell(8, 5).Range.Text = "19+43=62"
$t.Cell(9, 1).Range.Text = "91-34=57"
$t.Cell(9, 2).Range.Text = "7+39=46"
$t.Cell(9, 3).Range.Text = "61-33=28"
$t.Cell(9, 4).Range.Text = "38+15=53"
$t.Cell(9, 5).Range.Text = "9+88=97"
$t.Cell(10, 1).Range.Text = "47+14=61"
$t.Cell(10, 2).Range.Text = "2+59=61"
$t.Cell(10, 3).Range.Text = "25+28=53"
$t.Cell(10, 4).Range.Text = "19+26=45"
$t.Cell(10, 5).Range.Text = "94-48=46"
$t.Cell(11, 1).Range.Text = "54-7=47"
$t.Cell(11, 2).Range.Text = "84-5=79"
$t.Cell(11, 3).Range.Text = "30-12=18"
$t.Cell(11, 4).Range.Text = "17+14=31"
$t.Cell(11, 5).Range.Text = "19+48=67"
$t.Cell(12, 1).Range.Text = "38+7=45"
$t.Cell(12, 2).Range.Text = "31-3=28"
$t.Cell(12, 3).Range.Text = "46-7=39"
$t.Cell(12, 4).Range.Text = "93-7=86"
$t.Cell(12, 5).Range.Text = "39+54=93"
$t.Cell(13, 1).Range.Text = "5+57=62"
$t.Cell(13, 2).Range.Text = "57+36=93"
$t.Cell(13, 3).Range.Text = "50-18=32"
$t.Cell(13, 4).Range.Text = "91-33=58"
$t.Cell(13, 5).Range.Text = "9+64=73"
$t.Cell(14, 1).Range.Text = "16+76=92"
$t.Cell(14, 2).Range.Text = "56+18=74"
$t.Cell(14, 3).Range.Text = "62-53=9"
$t.Cell(14, 4).Range.Text = "32-9=23"
$t.Cell(14, 5).Range.Text = "3+79=82"
$t.Cell(15, 1).Range.Text = "39+16=55"
$t.Cell(15, 2).Range.Text = "93-25=68"
$t.Cell(15, 3).Range.Text = "77+5=82"
$t.Cell(15, 4).Range.Text = "39+29=68"
$t.Cell(15, 5).Range.Text = "83-36=47"
$t.Cell(16, 1).Range.Text = "38+18=56"
$t.Cell(16, 2).Range.Text = "22+69=91"
$t.Cell(16, 3).Range.Text = "87+6=93"
$t.Cell(16, 4).Range.Text = "39+56=95"
$t.Cell(16, 5).Range.Text = "84-9=75"
$t.Cell(17, 1).Range.Text = "24+58=82"
$t.Cell(17, 2).Range.Text = "82-17=65"
$t.Cell(17, 3).Range.Text = "26+45=71"
$t.Cell(17, 4).Range.Text = "90-13=77"
$t.Cell(17, 5).Range.Text = "9+82=91"
$t.Cell(18, 1).Range.Text = "81-34=47"
$t.Cell(18, 2).Range.Text = "25+58=83"
$t.Cell(18, 3).Range.Text = "49+7=56"
$t.Cell(18, 4).Range.Text = "15+46=61"
$t.Cell(18, 5).Range.Text = "84-49=35"
$t.Cell(19, 1).Range.Text = "70-29=41"
$t.Cell(19, 2).Range.Text = "91-3=88"
$t.Cell(19, 3).Range.Text = "63-36=27"
$t.Cell(19, 4).Range.Text = "60-23=37"
$t.Cell(19, 5).Range.Text = "65-17=48"
$t.Cell(20, 1).Range.Text = "65+7=72"
$t.Cell(20, 2).Range.Text = "37+46=83"
$t.Cell(20, 3).Range.Text = "72-15=57"
$t.Cell(20, 4).Range.Text = "71-29=42"
$t.Cell(20, 5).Range.Text = "73-55=18"
